# Entrega do projeto para coleção de Inverno
# Adds the Winter-collection categories to the "Planilha1" sheet and
# applies an integer number format to the whole ID column (B).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The stray styled-but-empty cell D15 from the previous edit is no longer
# needed - remove it (and its style) entirely.
$ws.Range("D15").Clear()

# New "Winter collection" categories appended below the existing table.
# Rows are written in the same (slightly out-of-order) sequence the
# original author used so new entries land in the same relative spot.
$ws.Range("A18").Value = "JERSEY"
$ws.Range("B18").Value = 33

$ws.Range("A19").Value = "CAMISA"
$ws.Range("B19").Value = 34

$ws.Range("A20").Value = "REGATA"
$ws.Range("B20").Value = 35

$ws.Range("A22").Value = "SAIA"
$ws.Range("B22").Value = 37

$ws.Range("A23").Value = "VESTIDO"
$ws.Range("B23").Value = 38

$ws.Range("A24").Value = "CUECA"
$ws.Range("B24").Value = 39

$ws.Range("A21").Value = "TOP"
$ws.Range("B21").Value = 36

# Format the whole ID column (B) as a plain integer.
$ws.Range("B1:B24").NumberFormat = "0"

# Iterative calculation (small MaxChange) as set in the updated workbook.
$excel.Iteration = $true
$excel.MaxIterations = 100
$excel.MaxChange = 0.0001

# Move/restore the active selection to M1, matching the saved view state.
[void]$ws.Range("M1").Select()
